# "Add more change on document"
# Adds a new "GSEID" variable row (with description) to the StoredVariableList
# sheet, then extends the same column A/B formatting down through row 42
# (mirrors the original author dragging the fill handle / pasting formats
# down well past the last data row), and finally moves the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data row (row 3): GSEID variable -------------------------------
$ws.Range("A3").Value = "GSEID"
$ws.Range("B3").Value = "sting"
$ws.Range("C3").Value = 'Store GSEID which used as data. Only useful when DataSourse == "Public"'

# Row 3's A/B cells pick up the same look as row 2 (Consolas styling).
$ws.Range("A2:B2").Copy()
$ws.Range("A3:B3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# --- Extend the A/B formatting further down (rows 4-38) -----------------
$ws.Range("A4:B38").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# --- Rows 39-42 only got column B formatted (no A cell) -----------------
$ws.Range("B2").Copy()
$ws.Range("B39:B42").PasteSpecial(-4122) | Out-Null # xlPasteFormats

$excel.CutCopyMode = 0

# --- Move the active selection, matching the saved view state -----------
$ws.Range("E14").Select() | Out-Null
